# Applies the "Completed UI of log, config screen" edit:
#  1. Typography sheet: row 6 (Mode) becomes the LogRecord entry (bigger size),
#     and a brand-new row 11 (LogColumn) is added using the same font/wildcard
#     settings as the other UI text entries.
#  2. Translation sheet: the old "UNLOCK LOG" heading row (row 14) is removed
#     (all rows below shift up by one), and the translation rows that used to
#     point at the generic "Mode" typography now point at the new "LogColumn"
#     or "LogRecord" typography entries.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Typography sheet
# ---------------------------------------------------------------------------
$typo = $wb.Worksheets.Item("Typography")

# Row 6 ("Mode") -> LogRecord, size 24 -> 26
$typo.Range("B6").Value = "LogRecord"
$typo.Range("D6").Value = 26

# New row 11: LogColumn entry (same font/wildcard settings as rows 8-10)
$typo.Range("B11").Value = "LogColumn"
$typo.Range("C11").Value = "KohinoorBangla.ttf"
$typo.Range("D11").Value = 26
$typo.Range("E11").Value = 4
$typo.Range("F11").Value = "?"
$typo.Range("G11").Value = "!`u{201D}`"#*%&()'`$+-@_, .:;?/~±×÷•º``´{}©£€^®¥_=[]¡¢|\¿><"
$typo.Range("H11").Value = "a-z,A-Z,0-9,0x0020-0x0060"
$typo.Range("I11").Value = ""

# ---------------------------------------------------------------------------
# 2. Translation sheet
# ---------------------------------------------------------------------------
$trans = $wb.Worksheets.Item("Translation")

# Remove the old "UNLOCK LOG" heading row - everything below shifts up by one.
$trans.Rows("14").Delete()

# Rows that used to reference the generic "Mode" typography now reference the
# new LogColumn / LogRecord typography entries (row numbers below are the
# *post-shift* row numbers).
$logColumnRows = @(6, 7, 15, 16, 19)
$logRecordRows = @(9, 10, 17, 18, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37)

foreach ($r in $logColumnRows) {
    $trans.Range("C$r").Value = "LogColumn"
}

foreach ($r in $logRecordRows) {
    $trans.Range("C$r").Value = "LogRecord"
}
